# "cambio inicio de sesion" - refresh the sample login/credentials data on
# the "Datos" sheet:
#   - username cell (A2) goes from a phone number to a login email address,
#     and becomes a mailto: hyperlink
#   - password cell (B2) gets a brand new password, in plain (unformatted)
#     text
#   - the stray duplicate phone number in C2 is cleared (its text format is
#     kept)
#   - the stray formatted cell left over at C6 is removed entirely
#   - a classification footer is added to the page setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- credentials rows -------------------------------------------------------

# A2: old phone number -> new login email, turned into a mailto: hyperlink
$ws.Range("A2").Value = "lorena.rodriguez@claro.com.co"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:lorena.rodriguez@claro.com.co") | Out-Null

# B2: old password -> new password, as plain text (no special number format)
$ws.Range("B2").Value = "Pruebas2026*"
$ws.Range("B2").ClearFormats() | Out-Null

# C2: used to hold the old phone number again; now left blank (keeps its
# existing text number-format)
$ws.Range("C2").ClearContents() | Out-Null

# stray leftover formatted cell is removed entirely (also shrinks the used
# range back down to row 4)
$ws.Range("C6").Clear() | Out-Null

# --- page setup / footer -----------------------------------------------------

$ws.PageSetup.CenterFooter = "_x000D_&1#&`"Aptos`"&6&K000000 Clasificaci" + [char]0x00F3 + "n: Uso Interno. Documento Claro Colombia"

# leave the selection where the file was last saved
$ws.Range("G4").Select() | Out-Null
